$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-02-25 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-26 Monday", 2)
$d.Content.Find.Execute("63÷5=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "67÷7=9, 4", 2)
$d.Content.Find.Execute("78÷2=39, 0", $true, $false, $false, $false, $false, $true, 1, $false, "45÷3=15, 0", 2)
$d.Content.Find.Execute("37÷5=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "51÷2=25, 1", 2)
$d.Content.Find.Execute("37÷3=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "96÷8=12, 0", 2)
$d.Content.Find.Execute("59÷8=7, 3", $true, $false, $false, $false, $false, $true, 1, $false, "17÷4=4, 1", 2)
$d.Content.Find.Execute("70÷5=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "70÷9=7, 7", 2)
$d.Content.Find.Execute("57÷9=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "35÷6=5, 5", 2)
$d.Content.Find.Execute("73÷7=10, 3", $true, $false, $false, $false, $false, $true, 1, $false, "83÷3=27, 2", 2)
$d.Content.Find.Execute("10÷5=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "56÷5=11, 1", 2)
$d.Content.Find.Execute("99÷9=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "32÷7=4, 4", 2)
$d.Content.Find.Execute("48÷4=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "32÷8=4, 0", 2)
$d.Content.Find.Execute("22÷7=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "47÷8=5, 7", 2)
$d.Content.Find.Execute("59÷9=6, 5", $true, $false, $false, $false, $false, $true, 1, $false, "62÷5=12, 2", 2)
$d.Content.Find.Execute("50÷9=5, 5", $true, $false, $false, $false, $false, $true, 1, $false, "96÷2=48, 0", 2)
$d.Content.Find.Execute("95÷3=31, 2", $true, $false, $false, $false, $false, $true, 1, $false, "90÷2=45, 0", 2)
$d.Content.Find.Execute("25÷8=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "35÷2=17, 1", 2)
$d.Content.Find.Execute("28÷3=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "83÷5=16, 3", 2)
$d.Content.Find.Execute("78÷7=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "75÷5=15, 0", 2)
$d.Content.Find.Execute("88÷7=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "64÷6=10, 4", 2)
$d.Content.Find.Execute("65÷7=9, 2", $true, $false, $false, $false, $false, $true, 1, $false, "31÷3=10, 1", 2)
$d.Content.Find.Execute("66÷7=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "57÷3=19, 0", 2)
$d.Content.Find.Execute("45÷7=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "36÷7=5, 1", 2)
$d.Content.Find.Execute("99÷5=19, 4", $true, $false, $false, $false, $false, $true, 1, $false, "49÷8=6, 1", 2)
$d.Content.Find.Execute("26÷4=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "18÷5=3, 3", 2)
$d.Content.Find.Execute("87÷7=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "60÷4=15, 0", 2)
